# Week 13 logging for the Football Team Players Data workbook.
# Updates weekly cumulative stat totals on both sheets and appends the
# newly-active players (W.Smallwood, L.Thomas) who saw their first
# touches this week.

$wb = $excel.ActiveWorkbook
$rushing = $wb.Worksheets.Item("Rushing")
$receiving = $wb.Worksheets.Item("Receiving")

# ---------------------------------------------------------------------
# Rushing: bump season totals for players who rushed this week
# ---------------------------------------------------------------------
$rushing.Range("E2").Value = 13
$rushing.Range("F2").Value = 9

$rushing.Range("C3").Value = 116
$rushing.Range("D3").Value = 76
$rushing.Range("E3").Value = 28
$rushing.Range("F3").Value = 42

$rushing.Range("D5").Value = 10

$rushing.Range("C7").Value = 3
$rushing.Range("E7").Value = 1

# New rushers this week: copy the formatting used by the rest of column A
# (bold/centered/bordered) onto the two new rows, then fill in the data.
$rushing.Range("A10").Copy()
$rushing.Range("A11").PasteSpecial(-4122)
$rushing.Range("A12").PasteSpecial(-4122)

$rushing.Range("A11").Value = 9
$rushing.Range("B11").Value = "W.Smallwood"
$rushing.Range("C11").Value = 0
$rushing.Range("D11").Value = 0
$rushing.Range("E11").Value = 1
$rushing.Range("F11").Value = 0

$rushing.Range("A12").Value = 10
$rushing.Range("B12").Value = "L.Thomas"
$rushing.Range("C12").Value = 0
$rushing.Range("D12").Value = 1
$rushing.Range("E12").Value = 0
$rushing.Range("F12").Value = 0

# ---------------------------------------------------------------------
# Receiving: bump season totals for players who caught passes this week
# ---------------------------------------------------------------------
$receiving.Range("C2").Value = 37
$receiving.Range("D2").Value = 32
$receiving.Range("G2").Value = 4
$receiving.Range("H2").Value = 3

$receiving.Range("C5").Value = 79
$receiving.Range("D5").Value = 56
$receiving.Range("E5").Value = 38

$receiving.Range("C6").Value = 7
$receiving.Range("D6").Value = 6

$receiving.Range("C7").Value = 34
$receiving.Range("D7").Value = 28

$receiving.Range("C8").Value = 14
$receiving.Range("D8").Value = 8
$receiving.Range("G8").Value = 3
$receiving.Range("H8").Value = 2

$receiving.Range("C10").Value = 22
$receiving.Range("D10").Value = 15

# New receivers (W.Smallwood, L.Thomas) are logged right after D.Carter,
# pushing J.Bates / R.Seals-Jones down two rows. Fix up formatting on the
# trailing rows first (row 14 reuses a stray leftover-formatted cell, row
# 15 is brand new), then write all the row contents.
$receiving.Range("A13").Copy()
$receiving.Range("A14").PasteSpecial(-4122)
$receiving.Range("A15").PasteSpecial(-4122)

$receiving.Range("B12").Value = "W.Smallwood"
$receiving.Range("C12").Value = 2
$receiving.Range("D12").Value = 2
$receiving.Range("E12").Value = 0
$receiving.Range("F12").Value = 0
$receiving.Range("G12").Value = 0
$receiving.Range("H12").Value = 0

$receiving.Range("B13").Value = "L.Thomas"
$receiving.Range("C13").Value = 2
$receiving.Range("D13").Value = 2
$receiving.Range("E13").Value = 3
$receiving.Range("F13").Value = 1
$receiving.Range("G13").Value = 1
$receiving.Range("H13").Value = 1

$receiving.Range("A14").Value = 12
$receiving.Range("B14").Value = "J.Bates"
$receiving.Range("C14").Value = 12
$receiving.Range("D14").Value = 11
$receiving.Range("E14").Value = 1
$receiving.Range("F14").Value = 1
$receiving.Range("G14").Value = 0
$receiving.Range("H14").Value = 0

$receiving.Range("A15").Value = 13
$receiving.Range("B15").Value = "R.Seals-Jones"
$receiving.Range("C15").Value = 20
$receiving.Range("D15").Value = 17
$receiving.Range("E15").Value = 2
$receiving.Range("F15").Value = 1
$receiving.Range("G15").Value = 7
$receiving.Range("H15").Value = 6

# The active/selected sheet moves from Receiving to Rushing in this revision.
$rushing.Activate()
